$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.692432045936584
$ws.Range("B1").Value = 2.070842742919922
$ws.Range("C1").Value = 5.321286678314209
$ws.Range("D1").Value = 1.345447897911072
$ws.Range("E1").Value = 0.750275194644928
